$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.980.94"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.676.83"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'215.22"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").Value = "'0.517"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.0619"
$ws.Range("D10").Value = "'20.30"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").Value = "1.912.81"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").Value = "1.659.56"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "26.981.05"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'237.21"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").Value = "'8.09"
$ws.Range("E19").Value = "  +4.98%  "
$ws.Range("D20").Value = "0.0₃0734"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  -0.71%  "
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("D25").Value = "'145.94"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("D26").Value = "'7.23"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").Value = "'16.11"
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("D30").Value = "'0.0499"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").Value = "1.480.46"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("E35").Value = "  +4.66%  "
$ws.Range("D36").Value = "'2.41"
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").Value = "'0.583"
$ws.Range("E37").Value = "  +1.72%  "
$ws.Range("E38").Value = "  +2.49%  "
$ws.Range("E39").Value = "  +0.70%  "
$ws.Range("E40").Value = "  -3.09%  "
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  +2.04%  "
$ws.Range("D44").Value = "'67.51"
$ws.Range("E44").Value = "  +1.78%  "
$ws.Range("D45").Value = "1.820.17"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("E46").Value = "  +0.78%  "
$ws.Range("D47").Value = "'90.50"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").Value = "0.0₆0107"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("D51").Value = "'0.0508"
$ws.Range("E51").Value = "  +0.15%  "
